$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.Value = "'26.250.23"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.Value = "'  -0.39%  "
$c.Style = "Normal"

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.Value = "'1.659.29"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.Value = "'  -0.51%  "
$c.Style = "Normal"

# Row 4
$c = $ws.Cells.Item(4, 5)
$c.Value = "'  -0.70%  "
$c.Style = "Normal"

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.Value = "'219.02"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.Value = "'  -0.26%  "
$c.Style = "Normal"

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.Value = "'0.5242"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.Value = "'  -1.95%  "
$c.Style = "Normal"

# Row 7
$c = $ws.Cells.Item(7, 5)
$c.Value = "'  -0.67%  "
$c.Style = "Normal"

# Row 8
$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.2647"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.Value = "'  -0.50%  "
$c.Style = "Normal"

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.06318"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.Value = "'  -1.15%  "
$c.Style = "Normal"

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.Value = "'20.70"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.Value = "'  -0.88%  "
$c.Style = "Normal"

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.07789"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.Value = "'  -0.69%  "
$c.Style = "Normal"

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.Value = "'4.513"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.Value = "'  -1.06%  "
$c.Style = "Normal"

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.Value = "'1.639.03"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.Value = "'  -1.75%  "
$c.Style = "Normal"

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.Value = "'1.888.96"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.Value = "'  -0.43%  "
$c.Style = "Normal"

# Row 15
$c = $ws.Cells.Item(15, 5)
$c.Value = "'  +1.66%  "
$c.Style = "Normal"

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.Value = "'0.0₅8065"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.Value = "'  -1.49%  "
$c.Style = "Normal"

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.Value = "'65.25"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.Value = "'  -1.08%  "
$c.Style = "Normal"

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.Value = "'26.245.69"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.Value = "'  -0.53%  "
$c.Style = "Normal"

# Row 19
$c = $ws.Cells.Item(19, 5)
$c.Value = "'  -0.65%  "
$c.Style = "Normal"

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.Value = "'4.724"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.Value = "'  +1.14%  "
$c.Style = "Normal"

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.Value = "'194.62"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.Value = "'  -0.11%  "
$c.Style = "Normal"

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.Value = "'10.24"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.Value = "'  -0.16%  "
$c.Style = "Normal"

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.Value = "'6.032"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.Value = "'  -0.06%  "
$c.Style = "Normal"

# Row 24
$c = $ws.Cells.Item(24, 5)
$c.Value = "'  -0.69%  "
$c.Style = "Normal"

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.Value = "'145.35"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.Value = "'  -0.76%  "
$c.Style = "Normal"

# Row 26
$c = $ws.Cells.Item(26, 5)
$c.Value = "'  -1.12%  "
$c.Style = "Normal"

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.Value = "'7.228"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.Value = "'  -0.13%  "
$c.Style = "Normal"

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.Value = "'16.08"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.Value = "'  -0.39%  "
$c.Style = "Normal"

# Row 29
$c = $ws.Cells.Item(29, 5)
$c.Value = "'  -0.20%  "
$c.Style = "Normal"

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.Value = "'0.05647"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.Value = "'  -3.61%  "
$c.Style = "Normal"

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.Value = "'1.277"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.Value = "'  -0.72%  "
$c.Style = "Normal"

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.Value = "'3.491"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.Value = "'  -2.44%  "
$c.Style = "Normal"

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.Value = "'3.372"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.Value = "'  +2.55%  "
$c.Style = "Normal"

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.Value = "'1.603"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.Value = "'  -0.59%  "
$c.Style = "Normal"

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.Value = "'2.811"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.Value = "'  -0.86%  "
$c.Style = "Normal"

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.Value = "'0.9450"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.Value = "'  -2.64%  "
$c.Style = "Normal"

# Row 37
$c = $ws.Cells.Item(37, 5)
$c.Value = "'  -0.80%  "
$c.Style = "Normal"

# Row 38
$c = $ws.Cells.Item(38, 4)
$c.Value = "'0.5756"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.Value = "'  -1.15%  "
$c.Style = "Normal"

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.Value = "'0.01607"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.Value = "'  +0.00%  "
$c.Style = "Normal"

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.Value = "'5.989"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.Value = "'  +2.24%  "
$c.Style = "Normal"

# Row 41
$c = $ws.Cells.Item(41, 2)
$c.Value = "'mCoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 3)
$c.Value = "'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.Value = "'2.568"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.Value = "'  -0.21%  "
$c.Style = "Normal"

# Row 42
$c = $ws.Cells.Item(42, 2)
$c.Value = "'Maker"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 3)
$c.Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 4)
$c.Value = "'1.050.43"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.Value = "'  -2.38%  "
$c.Style = "Normal"

# Row 43
$c = $ws.Cells.Item(43, 4)
$c.Value = "'0.8476"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.Value = "'  -1.91%  "
$c.Style = "Normal"

# Row 44
$c = $ws.Cells.Item(44, 5)
$c.Value = "'  -0.71%  "
$c.Style = "Normal"

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.Value = "'102.88"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.Value = "'  -1.45%  "
$c.Style = "Normal"

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.Value = "'1.798.90"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.Value = "'  -0.44%  "
$c.Style = "Normal"

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.Value = "'58.34"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.Value = "'  +0.62%  "
$c.Style = "Normal"

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.Value = "'0.0₈105"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.Value = "'  +0.56%  "
$c.Style = "Normal"

# Row 49
$c = $ws.Cells.Item(49, 5)
$c.Value = "'  -1.03%  "
$c.Style = "Normal"

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.Value = "'0.05322"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.Value = "'  +2.97%  "
$c.Style = "Normal"

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.Value = "'0.4355"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.Value = "'  -0.88%  "
$c.Style = "Normal"

